# RPA datasets push 2024-06-01
# Insert a new IPO row ("아이빔테크놀로지") at the top of the data table on the
# "02_38커뮤니케이션(최근일자기준)" sheet, shifting the existing rows down,
# and drop the last (oldest) row so the table keeps the same 20 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("02_38커뮤니케이션(최근일자기준)")

# Insert a new row above row 2 to make room for the new entry, pushing all
# the other data rows down by one.
$ws.Rows.Item(2).Insert()

# New top row: 아이빔테크놀로지
$ws.Range("A2").Value = "아이빔테크놀로지"
$ws.Range("B2").Value = "2024.07.15~07.19"
$ws.Range("C2").Value = "7,300~8,500"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = 16308
$ws.Range("F2").Value = "삼성증권"

# Drop the now-trailing 22nd row (old last entry, 씨어스테크놀로지) so the
# sheet keeps its original extent of 21 rows (20 data rows).
$ws.Rows.Item(22).Delete()

# The former last row (now row 21) keeps its company/date/price-range/offer
# columns from the prior row 20, but picks up new demand-amount and
# underwriter values.
$ws.Range("E21").Value = 13650
$ws.Range("F21").Value = "NH투자증권"
